# "added spec details view on project details page"
#
# Materials sheet gets two new rows of plywood spec data:
#   - a new "Select Maple" / "3/4 A1 Select Maple Plywood" row inserted
#     right after the existing "PG Maple" row (becomes new row 7)
#   - a new "Prefinished Maple" / "3/4 Prefinished Maple Plywood" row
#     inserted right before the last "Drawer Material" row (becomes new
#     row 21, pushing the old last row down to row 22)
# Selection/active-sheet state also moves around: Materials becomes the
# active tab/sheet (it was Drawers before), and a couple of other sheets
# simply have their remembered selection moved.

$wb = $excel.ActiveWorkbook

$wsMaterials = $wb.Worksheets.Item("Materials")

# --- Insert the two new rows -------------------------------------------
# Insert before current row 7 (shifts old rows 7-20 down to 8-21).
$wsMaterials.Rows.Item(7).Insert() | Out-Null
# Insert before (new) row 21 -- that's the old last row, which just moved
# from 20 to 21 because of the insert above. Shifts it down to row 22.
$wsMaterials.Rows.Item(21).Insert() | Out-Null

# --- Populate row 21 first so its new material string is created in the
# shared-string table ahead of row 7's (matches target string order) ----
$wsMaterials.Range("A21").Value = "Prefinished Maple"
$wsMaterials.Range("B21").Value = "3/4 Prefinished Maple Plywood"
$wsMaterials.Range("C21").Value = 0.75
$wsMaterials.Range("D21").Value = 48
$wsMaterials.Range("E21").Value = 96
$wsMaterials.Range("F21").Value = 40
$wsMaterials.Range("G21").Value = 0.1
$wsMaterials.Range("H21").Value = 0.2
$wsMaterials.Range("A21").HorizontalAlignment = -4131

# --- Now populate row 7 -------------------------------------------------
$wsMaterials.Range("A7").Value = "Select Maple"
$wsMaterials.Range("B7").Value = "3/4 A1 Select Maple Plywood"
$wsMaterials.Range("C7").Value = 0.75
$wsMaterials.Range("D7").Value = 48
$wsMaterials.Range("E7").Value = 96
$wsMaterials.Range("F7").Value = 120
$wsMaterials.Range("G7").Value = 0.15
$wsMaterials.Range("H7").Value = 0.2

# --- Move the remembered selection on a couple of other sheets ---------
$wsSpecs = $wb.Worksheets.Item("Specs")
$wsSpecs.Range("D31").Select() | Out-Null

$wsCabinets = $wb.Worksheets.Item("Cabinets")
$wsCabinets.Range("A52").Select() | Out-Null

$wsDrawers = $wb.Worksheets.Item("Drawers")
$wsDrawers.Range("H27").Select() | Out-Null

# --- Materials becomes the active sheet/tab, selection on A27 ----------
$wsMaterials.Activate() | Out-Null
$wsMaterials.Range("A27").Select() | Out-Null
